$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column C (current C..J metrics shift to D..K)
$ws.Range("C1").EntireColumn.Insert()

# Header for new column
$ws.Range("C1").Value = "M_PL"

# New column values (row 2..6)
$ws.Range("C2").Value = 1008209699708
$ws.Range("C3").Value = 3140810
$ws.Range("C4").Value = 21277927825
$ws.Range("C5").Value = 353160988340
$ws.Range("C6").Value = 49527932043
